# Fruta / hortaliza, semanal
# Re-shuffle the weekly data rows (dates + their associated price data)
# for "Hortaliza, Vega Monumental Concepción - Espárragos" so each row
# carries the record that actually belongs to it. Rows 1 (header), 2 and
# 8 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44496
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("J3").Value = 550
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1773
$ws.Range("N3").Value = '$/paquete'
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1773

# Row 4
$ws.Range("D4").Value = 44510
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1400
$ws.Range("M4").Value = 1350
$ws.Range("N4").Value = '$/kilo'
$ws.Range("O4").Value = "Provincia de Linares"
$ws.Range("P4").Value = 1350

# Row 5
$ws.Range("D5").Value = 44526
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1500
$ws.Range("L5").Value = 1600
$ws.Range("M5").Value = 1550
$ws.Range("N5").Value = '$/kilo'
$ws.Range("O5").Value = "Provincia de Linares"
$ws.Range("P5").Value = 1550

# Row 6
$ws.Range("D6").Value = 44545
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 550
$ws.Range("K6").Value = 1700
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = 1755
$ws.Range("N6").Value = '$/kilo'
$ws.Range("O6").Value = "Provincia de Linares"
$ws.Range("P6").Value = 1755

# Row 7
$ws.Range("D7").Value = 44477
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 1400
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1460
$ws.Range("N7").Value = '$/kilo'
$ws.Range("O7").Value = "Provincia de Linares"
$ws.Range("P7").Value = 1460

# Row 9
$ws.Range("D9").Value = 44468
$ws.Range("H9").Value = "Verde"
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1920
$ws.Range("N9").Value = '$/kilo'
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1920

# Row 10
$ws.Range("D10").Value = 44511
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1400
$ws.Range("M10").Value = 1350
$ws.Range("N10").Value = '$/kilo'
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1350

# Row 11
$ws.Range("D11").Value = 44524
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1600
$ws.Range("M11").Value = 1550
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = "Provincia de Talca"
$ws.Range("P11").Value = 1550
